# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E23) on Hoja1 currently lists the mora
# periods in descending order (2304 .. 2209). The updated account
# statement database lists them in ascending order (2209 .. 2304), and
# the "Valor Mora" amounts (F16:F23) must stay attached to their
# matching period (2209 -> 53120, 2304 -> 47808; the rest stay 53120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 16-23
$periodos = @("2209", "2210", "2211", "2212", "2301", "2302", "2303", "2304")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Valor Mora values follow their period: swap the values that were in the
# first and last rows of the table (period 2209 now on row 16, period 2304
# now on row 23).
$ws.Range("F16").Value = 53120
$ws.Range("F23").Value = 47808
